# Insert two new rows at the top of the 594-block and push the remaining
# data (previously rows 594-676) down by two rows to 596-678.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("594:595").Insert()

# New row 594
$ws.Range("A594").Value = 10
$ws.Range("B594").Value = "Vega Modelo de Temuco"
$ws.Range("C594").Value = "La Araucanía"
$ws.Range("D594").Value = 45131
$ws.Range("E594").Value = 9
$ws.Range("F594").Value = 100112008
$ws.Range("G594").Value = "Coliflor"
$ws.Range("H594").Value = "Sin especificar"
$ws.Range("I594").Value = "Primera"
$ws.Range("J594").Value = 1600
$ws.Range("K594").Value = 1200
$ws.Range("L594").Value = 1200
$ws.Range("M594").Value = 1200
$ws.Range("N594").Value = "$/unidad"
$ws.Range("O594").Value = "Región de O'Higgins"
$ws.Range("P594").Value = 1200
$ws.Range("Q594").Value = 1
$ws.Range("R594").Value = "Hortaliza"

# New row 595
$ws.Range("A595").Value = 10
$ws.Range("B595").Value = "Vega Modelo de Temuco"
$ws.Range("C595").Value = "La Araucanía"
$ws.Range("D595").Value = 45131
$ws.Range("E595").Value = 9
$ws.Range("F595").Value = 100112008
$ws.Range("G595").Value = "Coliflor"
$ws.Range("H595").Value = "Sin especificar"
$ws.Range("I595").Value = "Primera"
$ws.Range("J595").Value = 2180
$ws.Range("K595").Value = 1300
$ws.Range("L595").Value = 1300
$ws.Range("M595").Value = 1300
$ws.Range("N595").Value = "$/unidad"
$ws.Range("O595").Value = "Región del Maule"
$ws.Range("P595").Value = 1300
$ws.Range("Q595").Value = 1
$ws.Range("R595").Value = "Hortaliza"

